# Mark task "Tugas 3" (column E) as completed ("ü" checkmark in Wingdings
# font) for the three students whose row previously had this cell blank:
#   - E4  (row 4,  "Alief Faza Rizqi Adi Jaya")
#   - E7  (row 7,  "Cindy Anastasya Kurniawan Oey")
#   - E19 (row 19, "Muhammad Farhan")
#
# This mirrors the formatting already used by every other checked cell in
# the sheet: value "ü" rendered with the Wingdings font (size 12), which
# displays as a checkmark glyph. Only the font changes; fill/border stay
# as they already were on these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("E4", "E7", "E19")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "ü"
    $cell.Font.Name = "Wingdings"
    $cell.Font.Size = 12
}
